$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above row 415 (Excel shifts the old rows 415..442 down to 416..443,
# inheriting formatting from the row above, same as a manual "Insert" in the UI).
$ws.Rows.Item(415).Insert()

# Populate the newly inserted row with the new Ciboulette price observation.
$ws.Cells.Item(415, 1).Value = 6
$ws.Cells.Item(415, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(415, 3).Value = "Metropolitana"
$ws.Cells.Item(415, 4).Value = 44714
$ws.Cells.Item(415, 5).Value = 13
$ws.Cells.Item(415, 6).Value = 100112039
$ws.Cells.Item(415, 7).Value = "Ciboulette"
$ws.Cells.Item(415, 8).Value = "Sin especificar"
$ws.Cells.Item(415, 9).Value = "Primera"
$ws.Cells.Item(415, 10).Value = 740
$ws.Cells.Item(415, 11).Value = 700
$ws.Cells.Item(415, 12).Value = 800
$ws.Cells.Item(415, 13).Value = 747
$ws.Cells.Item(415, 14).Value = "$/docena de atados"
$ws.Cells.Item(415, 15).Value = "Región Metropolitana"
$ws.Cells.Item(415, 16).Value = 249
$ws.Cells.Item(415, 17).Value = 3
$ws.Cells.Item(415, 18).Value = "Hortaliza"
